$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1/IF in row 1, matching the style of the existing headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New data values for columns I and J
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 5

$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 8

$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 8
